# "Lagaði áætlun" - fix up the time-tracking plan on the diary sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Vika 6 (week 6, rows 26-32): log hours for "Rannsóknir" and "Prófanir"
$ws.Range("I26").Value = 2
$ws.Range("I29").Value = 3

# Vika 7 (week 7, rows 36-42): log hours for "Prófanir" and "Frágangur"
$ws.Range("C39").Value = 3
$ws.Range("C41").Value = 1

# Leave the view scrolled down / selection parked on the last edited cell,
# matching where the author was working when the plan was saved.
$excel.ActiveWindow.ScrollRow = 9
$ws.Range("D41").Select()
